$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 96.320746
$ws.Range("H2").Value = 288.962238
$ws.Range("I2").Value = 0.3809824610908788
$ws.Range("J2").Value = 0.3809824610908788
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 150.1098686666667
$ws.Range("N2").Value = 450.329606
$ws.Range("O2").Value = 0.7276622610660995
$ws.Range("P2").Value = 0.7276622610660997
$ws.Range("Q2").Value = 14458.69453193536
$ws.Range("R2").Value = 130128.2507874182
$ws.Range("S2").Value = 0.2772265590639162
$ws.Range("T2").Value = 0.2772265590639162

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 96.320746
$ws.Range("H3").Value = 288.962238
$ws.Range("I3").Value = 0.3809824610908788
$ws.Range("J3").Value = 0.3809824610908788
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 42.32476666666667
$ws.Range("N3").Value = 126.9743
$ws.Range("O3").Value = 0.2051706239258123
$ws.Range("P3").Value = 0.2051706239258124
$ws.Range("Q3").Value = 4076.753099609267
$ws.Range("R3").Value = 36690.7778964834
$ws.Range("S3").Value = 0.07816640924680712
$ws.Range("T3").Value = 0.07816640924680714

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 96.320746
$ws.Range("H4").Value = 288.962238
$ws.Range("I4").Value = 0.3809824610908788
$ws.Range("J4").Value = 0.3809824610908788
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.940565666666667
$ws.Range("N4").Value = 14.821697
$ws.Range("O4").Value = 0.02394954586187395
$ws.Range("P4").Value = 0.02394954586187395
$ws.Range("Q4").Value = 475.8789706753207
$ws.Range("R4").Value = 4282.910736077886
$ws.Range("S4").Value = 0.009124356924465608
$ws.Range("T4").Value = 0.00912435692446561

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 96.320746
$ws.Range("H5").Value = 288.962238
$ws.Range("I5").Value = 0.3809824610908788
$ws.Range("J5").Value = 0.3809824610908788
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.915377333333334
$ws.Range("N5").Value = 26.746132
$ws.Range("O5").Value = 0.04321756914621411
$ws.Range("P5").Value = 0.04321756914621412
$ws.Range("Q5").Value = 858.7357956181573
$ws.Range("R5").Value = 7728.622160563416
$ws.Range("S5").Value = 0.01646513585568988
$ws.Range("T5").Value = 0.01646513585568989

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.46467533333333
$ws.Range("H6").Value = 55.394026
$ws.Range("I6").Value = 0.07303429161291354
$ws.Range("J6").Value = 0.07303429161291354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 150.1098686666667
$ws.Range("N6").Value = 450.329606
$ws.Range("O6").Value = 0.7276622610660995
$ws.Range("P6").Value = 0.7276622610660997
$ws.Range("Q6").Value = 2771.729989259306
$ws.Range("R6").Value = 24945.56990333376
$ws.Range("S6").Value = 0.05314429777041354
$ws.Range("T6").Value = 0.05314429777041355

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.46467533333333
$ws.Range("H7").Value = 55.394026
$ws.Range("I7").Value = 0.07303429161291354
$ws.Range("J7").Value = 0.07303429161291354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 42.32476666666667
$ws.Range("N7").Value = 126.9743
$ws.Range("O7").Value = 0.2051706239258123
$ws.Range("P7").Value = 0.2051706239258124
$ws.Range("Q7").Value = 781.5130750590889
$ws.Range("R7").Value = 7033.617675531799
$ws.Range("S7").Value = 0.01498449117820119
$ws.Range("T7").Value = 0.0149844911782012

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.46467533333333
$ws.Range("H8").Value = 55.394026
$ws.Range("I8").Value = 0.07303429161291354
$ws.Range("J8").Value = 0.07303429161291354
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.940565666666667
$ws.Range("N8").Value = 14.821697
$ws.Range("O8").Value = 0.02394954586187395
$ws.Range("P8").Value = 0.02394954586187395
$ws.Range("Q8").Value = 91.22594099801356
$ws.Range("R8").Value = 821.0334689821219
$ws.Range("S8").Value = 0.001749138116472949
$ws.Range("T8").Value = 0.001749138116472949

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.46467533333333
$ws.Range("H9").Value = 55.394026
$ws.Range("I9").Value = 0.07303429161291354
$ws.Range("J9").Value = 0.07303429161291354
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.915377333333334
$ws.Range("N9").Value = 26.746132
$ws.Range("O9").Value = 0.04321756914621411
$ws.Range("P9").Value = 0.04321756914621412
$ws.Range("Q9").Value = 164.6195479341591
$ws.Range("R9").Value = 1481.575931407432
$ws.Range("S9").Value = 0.003156364547825856
$ws.Range("T9").Value = 0.003156364547825857

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 122.909391
$ws.Range("H10").Value = 368.728173
$ws.Range("I10").Value = 0.4861499128584522
$ws.Range("J10").Value = 0.4861499128584522
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 150.1098686666667
$ws.Range("N10").Value = 450.329606
$ws.Range("O10").Value = 0.7276622610660995
$ws.Range("P10").Value = 0.7276622610660997
$ws.Range("Q10").Value = 18449.91254090998
$ws.Range("R10").Value = 166049.2128681898
$ws.Range("S10").Value = 0.3537529448076686
$ws.Range("T10").Value = 0.3537529448076686

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 122.909391
$ws.Range("H11").Value = 368.728173
$ws.Range("I11").Value = 0.4861499128584522
$ws.Range("J11").Value = 0.4861499128584522
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 42.32476666666667
$ws.Range("N11").Value = 126.9743
$ws.Range("O11").Value = 0.2051706239258123
$ws.Range("P11").Value = 0.2051706239258124
$ws.Range("Q11").Value = 5202.1112952171
$ws.Range("R11").Value = 46819.00165695389
$ws.Range("S11").Value = 0.09974368094264793
$ws.Range("T11").Value = 0.09974368094264795

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 122.909391
$ws.Range("H12").Value = 368.728173
$ws.Range("I12").Value = 0.4861499128584522
$ws.Range("J12").Value = 0.4861499128584522
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.940565666666667
$ws.Range("N12").Value = 14.821697
$ws.Range("O12").Value = 0.02394954586187395
$ws.Range("P12").Value = 0.02394954586187395
$ws.Range("Q12").Value = 607.241917285509
$ws.Range("R12").Value = 5465.177255569581
$ws.Range("S12").Value = 0.01164306963374952
$ws.Range("T12").Value = 0.01164306963374952

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 122.909391
$ws.Range("H13").Value = 368.728173
$ws.Range("I13").Value = 0.4861499128584522
$ws.Range("J13").Value = 0.4861499128584522
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.915377333333334
$ws.Range("N13").Value = 26.746132
$ws.Range("O13").Value = 0.04321756914621411
$ws.Range("P13").Value = 0.04321756914621412
$ws.Range("Q13").Value = 1095.783598575204
$ws.Range("R13").Value = 9862.052387176835
$ws.Range("S13").Value = 0.02101021747438612
$ws.Range("T13").Value = 0.02101021747438613

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 15.127183
$ws.Range("H14").Value = 45.381549
$ws.Range("I14").Value = 0.05983333443775553
$ws.Range("J14").Value = 0.05983333443775553
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 150.1098686666667
$ws.Range("N14").Value = 450.329606
$ws.Range("O14").Value = 0.7276622610660995
$ws.Range("P14").Value = 0.7276622610660997
$ws.Range("Q14").Value = 2270.739453426633
$ws.Range("R14").Value = 20436.65508083969
$ws.Range("S14").Value = 0.0435384594241013
$ws.Range("T14").Value = 0.04353845942410131

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 15.127183
$ws.Range("H15").Value = 45.381549
$ws.Range("I15").Value = 0.05983333443775553
$ws.Range("J15").Value = 0.05983333443775553
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 42.32476666666667
$ws.Range("N15").Value = 126.9743
$ws.Range("O15").Value = 0.2051706239258123
$ws.Range("P15").Value = 0.2051706239258124
$ws.Range("Q15").Value = 640.2544907989667
$ws.Range("R15").Value = 5762.2904171907
$ws.Range("S15").Value = 0.0122760425581561
$ws.Range("T15").Value = 0.0122760425581561

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 15.127183
$ws.Range("H16").Value = 45.381549
$ws.Range("I16").Value = 0.05983333443775553
$ws.Range("J16").Value = 0.05983333443775553
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 4.940565666666667
$ws.Range("N16").Value = 14.821697
$ws.Range("O16").Value = 0.02394954586187395
$ws.Range("P16").Value = 0.02394954586187395
$ws.Range("Q16").Value = 74.73684096318368
$ws.Range("R16").Value = 672.631568668653
$ws.Range("S16").Value = 0.001432981187185868
$ws.Range("T16").Value = 0.001432981187185868

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 15.127183
$ws.Range("H17").Value = 45.381549
$ws.Range("I17").Value = 0.05983333443775553
$ws.Range("J17").Value = 0.05983333443775553
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 8.915377333333334
$ws.Range("N17").Value = 26.746132
$ws.Range("O17").Value = 0.04321756914621411
$ws.Range("P17").Value = 0.04321756914621412
$ws.Range("Q17").Value = 134.8645444353853
$ws.Range("R17").Value = 1213.780899918468
$ws.Range("S17").Value = 0.002585851268312254
$ws.Range("T17").Value = 0.002585851268312254
